# Prefix the image filenames in column A ("Image") with "existence/",
# e.g. "existence_0.png" -> "existence/existence_0.png".
# Row 1 is the header ("Image") and is left untouched; data starts at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($current -and $current -notlike "existence/*") {
        $cell.Value2 = "existence/" + $current
    }
}
